$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting Time Period/Comments/Win-Lost to B/C/D
$ws.Range("A1").EntireColumn.Insert() | Out-Null

# New header cell for the inserted column, matching the style of the other headers
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$ws.Range("A1").Value = "Date"

# Clear out the old data rows (previously rows 3 and 4, now 3 and 4 after column insert)
$ws.Rows("3:4").Delete() | Out-Null

# Write the single consolidated data row into row 2
# Force text format on A2 so the date-like string is stored as text (not a date serial),
# then restore the default style so the cell keeps an unstyled/default appearance.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2024-08-30"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "3-4"
$ws.Range("C2").Value = "cdsf"
$ws.Range("D2").Value = "fsdf"
